$d = $word.ActiveDocument

# Locate the last paragraph in the document body (the "Interviewer" paragraph)
# and build an insertion point immediately after it (but still before sectPr).
$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range
$insertPoint = $d.Range($lastRange.End, $lastRange.End)

# WordprocessingML fragment for the new "TAGS:" heading paragraph plus the
# twelve tag-list paragraphs that follow it.
$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>TAGS:</w:t></w:r></w:p><w:p><w:r><w:t>1 Nepiaq, Enepiaq, Enpiaq -- Sod-house</w:t></w:r></w:p><w:p><w:r><w:t>1 Qanruyutet, Qaneryarat -- Traditional Wisdom, Wise Words</w:t></w:r></w:p><w:p><w:r><w:t>1 Aatailnguut -- Illegitimate Children</w:t></w:r></w:p><w:p><w:r><w:t>1 Tan''gaurluut Nasaurluut-llu Allakarluteng -- Boys and Girls Separate</w:t></w:r></w:p><w:p><w:r><w:t>1 Yugtun Qalarcaraq, Qaneryaraq -- Yup''ik Language</w:t></w:r></w:p><w:p><w:r><w:t>1 Agayumayaraq, Agayumaciq -- Christianity</w:t></w:r></w:p><w:p><w:r><w:t>1 Taangaq -- Alcohol</w:t></w:r></w:p><w:p><w:r><w:t>1 Ilaliuryaraq -- Socializing</w:t></w:r></w:p><w:p><w:r><w:t>1 Ayagyuat -- Youth</w:t></w:r></w:p><w:p><w:r><w:t>1 Maligtaquyaraq, Niisngayaraq -- Obedience</w:t></w:r></w:p><w:p><w:r><w:t>1 Ellminek Tuqutellerkaq -- Suicide</w:t></w:r></w:p><w:p><w:r><w:t>1 Piicak -- Prayer</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertPoint.InsertXML($xmlFrag)
